$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$reqText = "LOM3202 -  Circuitos Elétricos  (Requisito)`n"
$setText = "LOM3206 -  Eletrônica  (Indicação de Conjunto)`n"

$ws.Range("B24").Value = $setText
$ws.Range("C24").Value = $setText
$ws.Range("B25").Value = $reqText
$ws.Range("C25").Value = $reqText
